# Update "Forecast Comparison" sheet:
#  - Insert a new column B named "Week_Start_Date" holding the first day of
#    each forecast week (shifts ASIN..is_holiday_week one column to the right).
#  - Change the Week labels in column A from zero-padded "W01".."W16" to
#    "W1".."W16".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new, blank column before column B. This pushes the existing
# ASIN / MyForecast / ... / is_holiday_week columns from B:I to C:J.
$ws.Columns("B:B").Insert()

# Header for the newly inserted column.
$ws.Range("B1").Value = "Week_Start_Date"

# Week number labels (column A) and week-start dates (column B) for the
# 16 data rows. Dates are stored as plain text (e.g. "2025-01-05"), not as
# Excel date values, so force the cells to text format before assigning.
$weeks = @(
    @{Row = 2;  Label = "W1";  Date = "2025-01-05"},
    @{Row = 3;  Label = "W2";  Date = "2025-01-12"},
    @{Row = 4;  Label = "W3";  Date = "2025-01-19"},
    @{Row = 5;  Label = "W4";  Date = "2025-01-26"},
    @{Row = 6;  Label = "W5";  Date = "2025-02-02"},
    @{Row = 7;  Label = "W6";  Date = "2025-02-09"},
    @{Row = 8;  Label = "W7";  Date = "2025-02-16"},
    @{Row = 9;  Label = "W8";  Date = "2025-02-23"},
    @{Row = 10; Label = "W9";  Date = "2025-03-02"},
    @{Row = 11; Label = "W10"; Date = "2025-03-09"},
    @{Row = 12; Label = "W11"; Date = "2025-03-16"},
    @{Row = 13; Label = "W12"; Date = "2025-03-23"},
    @{Row = 14; Label = "W13"; Date = "2025-03-30"},
    @{Row = 15; Label = "W14"; Date = "2025-04-06"},
    @{Row = 16; Label = "W15"; Date = "2025-04-13"},
    @{Row = 17; Label = "W16"; Date = "2025-04-20"}
)

foreach ($week in $weeks) {
    $r = $week.Row
    $ws.Cells.Item($r, 1).Value = $week.Label

    $dateCell = $ws.Cells.Item($r, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $week.Date
}
